$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 300
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 300
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -185
$ws.Range("N5").Value = $null
$ws.Range("H9").Value = 5000045.5
$ws.Range("I9").Value = 5000045.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 5000045.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -4999876.5
$ws.Range("N9").Value = $null
$ws.Range("H33").Value = 401102.4
$ws.Range("I33").Value = 909249
$ws.Range("J33").Value = 1844.3572
$ws.Range("K33").Value = 909249
$ws.Range("L33").Value = 1844.3572
$ws.Range("M33").Value = -909020
$ws.Range("N33").Value = -2302.3572
$ws.Range("H74").Value = 3208.818
$ws.Range("I74").Value = 3208.818
$ws.Range("K74").Value = 3208.818
$ws.Range("M74").Value = -2272.818
$ws.Range("H76").Value = 4994.375
$ws.Range("I76").Value = 4993.5713
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 4993.5713
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -4678.5713
$ws.Range("N76").Value = -5630
$ws.Range("H77").Value = 3208.818
$ws.Range("I77").Value = 3208.818
$ws.Range("K77").Value = 16044.09
$ws.Range("M77").Value = -11364.09
$ws.Range("H79").Value = 4994.375
$ws.Range("I79").Value = 4993.5713
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 4993.5713
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -3901.5713
$ws.Range("N79").Value = -7184
$ws.Range("H86").Value = 8900
$ws.Range("J86").Value = 8900
$ws.Range("L86").Value = 8900
$ws.Range("N86").Value = -11146
$ws.Range("H89").Value = 8900
$ws.Range("J89").Value = 8900
$ws.Range("L89").Value = 44500
$ws.Range("N89").Value = -55732
$ws.Range("H100").Value = 3009.818
$ws.Range("I100").Value = 1684.2858
$ws.Range("K100").Value = 1684.2858
$ws.Range("M100").Value = -1143.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 6080
$ws.Range("I22").Value = 4144.2856
$ws.Range("K22").Value = 4144.2856
$ws.Range("M22").Value = -3845.2856
$ws.Range("H32").Value = 5221.3267
$ws.Range("I32").Value = 3681.889
$ws.Range("K32").Value = 3681.889
$ws.Range("M32").Value = -3394.889
$ws.Range("H33").Value = 8700
$ws.Range("I33").Value = 8700
$ws.Range("K33").Value = 8700
$ws.Range("M33").Value = -8371
$ws.Range("H45").Value = 9956.385
$ws.Range("I45").Value = 14678.875
$ws.Range("J45").Value = 2400.4
$ws.Range("K45").Value = 14678.875
$ws.Range("L45").Value = 2400.4
$ws.Range("M45").Value = -14301.875
$ws.Range("N45").Value = -3154.4
$ws.Range("H139").Value = 82000
$ws.Range("I139").Value = 59000
$ws.Range("J139").Value = 105000
$ws.Range("K139").Value = 59000
$ws.Range("L139").Value = 105000
$ws.Range("M139").Value = -53860
$ws.Range("N139").Value = -115280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13326.272
$ws.Range("I26").Value = 13326.272
$ws.Range("K26").Value = 13326.272
$ws.Range("M26").Value = -13034.272
$ws.Range("H86").Value = 8451.048000000001
$ws.Range("I86").Value = 13775.111
$ws.Range("K86").Value = 13775.111
$ws.Range("M86").Value = -12652.111
$ws.Range("H89").Value = 8451.048000000001
$ws.Range("I89").Value = 13775.111
$ws.Range("K89").Value = 68875.55500000001
$ws.Range("M89").Value = -63259.55500000001
$ws.Range("H94").Value = 1720.2858
$ws.Range("I94").Value = 1720.2858
$ws.Range("K94").Value = 1720.2858
$ws.Range("M94").Value = -1269.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 99999
$ws.Range("J70").Value = 99999
$ws.Range("L70").Value = 99999
$ws.Range("N70").Value = -100629
$ws.Range("H73").Value = 99999
$ws.Range("J73").Value = 99999
$ws.Range("L73").Value = 99999
$ws.Range("N73").Value = -102183
$ws.Range("H86").Value = 6122.7
$ws.Range("I86").Value = 5896.933
$ws.Range("J86").Value = 6800
$ws.Range("K86").Value = 5896.933
$ws.Range("L86").Value = 6800
$ws.Range("M86").Value = -4773.933
$ws.Range("N86").Value = -9046
$ws.Range("H89").Value = 6122.7
$ws.Range("I89").Value = 5896.933
$ws.Range("J89").Value = 6800
$ws.Range("K89").Value = 29484.665
$ws.Range("L89").Value = 34000
$ws.Range("M89").Value = -23868.665
$ws.Range("N89").Value = -45232
$ws.Range("H93").Value = 17299.889
$ws.Range("I93").Value = 15712.375
$ws.Range("K93").Value = 15712.375
$ws.Range("M93").Value = -13840.375
$ws.Range("H105").Value = 8570.571
$ws.Range("I105").Value = 6999.25
$ws.Range("K105").Value = 6999.25
$ws.Range("M105").Value = -5252.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 1156
$ws.Range("I70").Value = 1156
$ws.Range("K70").Value = 3468
$ws.Range("M70").Value = -3153
$ws.Range("H73").Value = 1156
$ws.Range("I73").Value = 1156
$ws.Range("K73").Value = 3468
$ws.Range("M73").Value = -2376
$ws.Range("H132").Value = 2999.7144
$ws.Range("J132").Value = 2999.7144
$ws.Range("L132").Value = 26997.4296
$ws.Range("N132").Value = -32057.4296
$ws.Range("H137").Value = 9248.429
$ws.Range("I137").Value = 7874
$ws.Range("J137").Value = 9571.823
$ws.Range("K137").Value = 23622
$ws.Range("L137").Value = 28715.469
$ws.Range("M137").Value = -18522
$ws.Range("N137").Value = -38915.469
$ws.Range("H140").Value = 2774.1538
$ws.Range("I140").Value = 2169.25
$ws.Range("K140").Value = 6507.75
$ws.Range("M140").Value = -1327.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 594009600
$ws.Range("J15").Value = 594009600
$ws.Range("L15").Value = 594009600
$ws.Range("N15").Value = -594010176
$ws.Range("H70").Value = 11452.059
$ws.Range("I70").Value = 11715.333
$ws.Range("K70").Value = 11715.333
$ws.Range("M70").Value = -11445.333
$ws.Range("H73").Value = 11452.059
$ws.Range("I73").Value = 11715.333
$ws.Range("K73").Value = 11715.333
$ws.Range("M73").Value = -10779.333
$ws.Range("H81").Value = 594009600
$ws.Range("J81").Value = 594009600
$ws.Range("L81").Value = 594009600
$ws.Range("N81").Value = -594011596
$ws.Range("H84").Value = 594009600
$ws.Range("J84").Value = 594009600
$ws.Range("L84").Value = 1782028800
$ws.Range("N84").Value = -1782038784
$ws.Range("H100").Value = 48131
$ws.Range("J100").Value = 48131
$ws.Range("L100").Value = 48131
$ws.Range("N100").Value = -50295
$ws.Range("H102").Value = 3657.8572
$ws.Range("I102").Value = 2522.2
$ws.Range("K102").Value = 2522.2
$ws.Range("M102").Value = -900.1999999999998
$ws.Range("H123").Value = 43959
$ws.Range("J123").Value = 43959
$ws.Range("L123").Value = 43959
$ws.Range("N123").Value = -48859
$ws.Range("H126").Value = 2788.4285
$ws.Range("I126").Value = 2086.5
$ws.Range("K126").Value = 6259.5
$ws.Range("M126").Value = -3789.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4195.533
$ws.Range("I16").Value = 906.7368
$ws.Range("K16").Value = 906.7368
$ws.Range("M16").Value = -736.7368
$ws.Range("H40").Value = 4445.725
$ws.Range("I40").Value = 3587.3044
$ws.Range("K40").Value = 3587.3044
$ws.Range("M40").Value = -3451.3044
$ws.Range("H100").Value = 5492.242
$ws.Range("I100").Value = 4077.6
$ws.Range("J100").Value = 9913
$ws.Range("K100").Value = 4077.6
$ws.Range("L100").Value = 9913
$ws.Range("M100").Value = -3536.6
$ws.Range("N100").Value = -10995
$ws.Range("H122").Value = 2990.5652
$ws.Range("I122").Value = 2831.4736
$ws.Range("K122").Value = 8494.4208
$ws.Range("M122").Value = -6044.4208

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H100").Value = 891
$ws.Range("I100").Value = 410.22223
$ws.Range("J100").Value = 2333.3333
$ws.Range("K100").Value = 820.44446
$ws.Range("L100").Value = 4666.6666
$ws.Range("M100").Value = -279.44446
$ws.Range("N100").Value = -5748.6666
$ws.Range("H122").Value = 4754.591
$ws.Range("I122").Value = 2161.7693
$ws.Range("K122").Value = 6485.3079
$ws.Range("M122").Value = -4035.3079
